$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UnitMod")

# Commit: "Changed 1% modifier on units vs units to 5%"
# The 0.01 "units vs units" self-matchup modifiers in the R (cata/hele) and
# S (cannon) columns for the Knight (row 21), Archer (row 22), Hoplite
# (row 26) and Wagon (row 29) source rows become 0.05. Every other cell
# touched by the published diff (rows 43/44/48, 62/63/67, 80/81/85/88,
# 99/100/104/105, 116/121/133/137/140) is a formula that derives from
# these eight cells, so updating the sources and letting Excel recalc
# reproduces the same cached results.
$ws.Range("R21").Value = 0.05
$ws.Range("S21").Value = 0.05

$ws.Range("R22").Value = 0.05
$ws.Range("S22").Value = 0.05

$ws.Range("R26").Value = 0.05
$ws.Range("S26").Value = 0.05

$ws.Range("R29").Value = 0.05
$ws.Range("S29").Value = 0.05

# The author's saved cursor position ended up on R29 (single cell) rather
# than the original D15:V30 block selection.
[void]$ws.Range("R29").Select()
